# Generate Report for Handoff
#
# Rows 2 and 3 of each sheet (Overview, zh-cn, de-de) swap identity:
#   - the "acc75e5c-...md" file now occupies row 2 (status stays
#     "Handed back: in sync with en-US")
#   - the "a1aea548-...md" file now occupies row 3 and its status flips to
#     "Ready for handoff"; its handoff datetime is refreshed to reflect the
#     new handoff.
# Hyperlink targets (r:id) stay pinned to their cell address; only the
# visible display text is swapped to match the new cell content.

$wb = $excel.ActiveWorkbook

$accFile = "acc75e5c-c5dd-42b1-93da-3bb4b02724d1.md"
$a1aFile = "a1aea548-e9aa-4f55-8bde-b8729672a656.md"
$statusHandedBack = "Handed back: in sync with en-US"
$statusReady = "Ready for handoff"
$include = "Include"

function Set-HyperlinkDisplay($ws, $addr, $text) {
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $addr) {
            $hl.TextToDisplay = $text
        }
    }
}

# ---------------------------------------------------------------
# Sheet "Overview": columns A (file), B (zh-cn status), C (de-de status)
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $accFile
$wsOverview.Range("B2").Value = $statusHandedBack
$wsOverview.Range("C2").Value = $statusHandedBack

$wsOverview.Range("A3").Value = $a1aFile
$wsOverview.Range("B3").Value = $statusReady
$wsOverview.Range("C3").Value = $statusReady

Set-HyperlinkDisplay $wsOverview '$A$2' $accFile
Set-HyperlinkDisplay $wsOverview '$A$3' $a1aFile

# ---------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$zhAccHandoff = "acc75e5c-c5dd-42b1-93da-3bb4b02724d1.b49e60cd653bb0e7730d8174e1922021f3dfb690.zh-cn.xlf"
$zhA1aHandoff = "a1aea548-e9aa-4f55-8bde-b8729672a656.646c5bde1652c0fd3119e81a69089df399eb537a.zh-cn.xlf"

# Row 2 -> acc75e5c (handed back, in sync)
$wsZh.Range("A2").Value = $accFile
$wsZh.Range("B2").Value = $statusHandedBack
$wsZh.Range("C2").Value = $zhAccHandoff
$wsZh.Range("D2").Value = "2016-03-07 04:38:09"
$wsZh.Range("E2").Value = $accFile
$wsZh.Range("F2").Value = $zhAccHandoff
$wsZh.Range("G2").Value = "2016-03-07 04:39:06"
$wsZh.Range("H2").Value = $include

# Row 3 -> a1aea548 (ready for a new handoff)
$wsZh.Range("A3").Value = $a1aFile
$wsZh.Range("B3").Value = $statusReady
$wsZh.Range("C3").Value = $zhA1aHandoff
$wsZh.Range("D3").Value = "2016-03-07 04:40:14"
$wsZh.Range("E3").Value = $a1aFile
$wsZh.Range("F3").Value = $zhA1aHandoff
$wsZh.Range("G3").Value = "2016-03-07 04:39:06"
$wsZh.Range("H3").Value = $include

Set-HyperlinkDisplay $wsZh '$A$2' $accFile
Set-HyperlinkDisplay $wsZh '$C$2' $zhAccHandoff
Set-HyperlinkDisplay $wsZh '$E$2' $accFile
Set-HyperlinkDisplay $wsZh '$F$2' $zhAccHandoff
Set-HyperlinkDisplay $wsZh '$A$3' $a1aFile
Set-HyperlinkDisplay $wsZh '$C$3' $zhA1aHandoff
Set-HyperlinkDisplay $wsZh '$E$3' $a1aFile
Set-HyperlinkDisplay $wsZh '$F$3' $zhA1aHandoff

# ---------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$deAccHandoff = "acc75e5c-c5dd-42b1-93da-3bb4b02724d1.b49e60cd653bb0e7730d8174e1922021f3dfb690.de-de.xlf"
$deA1aHandoff = "a1aea548-e9aa-4f55-8bde-b8729672a656.646c5bde1652c0fd3119e81a69089df399eb537a.de-de.xlf"

# Row 2 -> acc75e5c (handed back, in sync)
$wsDe.Range("A2").Value = $accFile
$wsDe.Range("B2").Value = $statusHandedBack
$wsDe.Range("C2").Value = $deAccHandoff
$wsDe.Range("D2").Value = "2016-03-07 04:38:23"
$wsDe.Range("E2").Value = $accFile
$wsDe.Range("F2").Value = $deAccHandoff
$wsDe.Range("G2").Value = "2016-03-07 04:39:25"
$wsDe.Range("H2").Value = $include

# Row 3 -> a1aea548 (ready for a new handoff)
$wsDe.Range("A3").Value = $a1aFile
$wsDe.Range("B3").Value = $statusReady
$wsDe.Range("C3").Value = $deA1aHandoff
$wsDe.Range("D3").Value = "2016-03-07 04:40:25"
$wsDe.Range("E3").Value = $a1aFile
$wsDe.Range("F3").Value = $deA1aHandoff
$wsDe.Range("G3").Value = "2016-03-07 04:39:25"
$wsDe.Range("H3").Value = $include

Set-HyperlinkDisplay $wsDe '$A$2' $accFile
Set-HyperlinkDisplay $wsDe '$C$2' $deAccHandoff
Set-HyperlinkDisplay $wsDe '$E$2' $accFile
Set-HyperlinkDisplay $wsDe '$F$2' $deAccHandoff
Set-HyperlinkDisplay $wsDe '$A$3' $a1aFile
Set-HyperlinkDisplay $wsDe '$C$3' $deA1aHandoff
Set-HyperlinkDisplay $wsDe '$E$3' $a1aFile
Set-HyperlinkDisplay $wsDe '$F$3' $deA1aHandoff
